# Atualizacao dos dados da bibi
# 1) Corrige o valor de faturamento do dia 4 de junho/2025 (linha 5)
# 2) Insere um novo registro para o dia 5 de junho/2025 (nova linha 6),
#    deslocando as linhas seguintes uma posicao para baixo.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrige o valor da linha 5 (dia 4, junho/2025)
$ws.Range("B5").Value = 36189.46

# Insere uma nova linha na posicao 6, empurrando os dados existentes para baixo
$ws.Rows.Item(6).Insert()

# Preenche os dados do novo dia (dia 5, junho/2025)
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = 21588.35
$ws.Cells.Item(6, 3).Value = 6
$ws.Cells.Item(6, 4).Value = 2025
$ws.Cells.Item(6, 5).Value = "06/2025"
